# Applies the coin price / 1h-volume refresh described in the commit
# "Updated cryptos list on Sun Jan 28 13:30:23 UTC 2024 with GitHub Actions".
# Rows 44/45 additionally swap their ranking order (VeChain <-> EnergySwap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.349.67"
$ws.Range("E2").Value = "  +1.33%  "

$ws.Range("D3").Value = "2.277.84"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'307.09"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").Value = "'98.04"
$ws.Range("E6").Value = "  +6.07%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.494"
$ws.Range("E9").Value = "  +2.27%  "

$ws.Range("D10").Value = "'35.78"
$ws.Range("E10").Value = "  +10.07%  "

$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D13").Value = "'6.69"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").Value = "2.627.36"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("D15").Value = "'14.44"
$ws.Range("E15").Value = "  +1.36%  "

$ws.Range("D16").Value = "2.280.19"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").Value = "'0.795"
$ws.Range("E17").Value = "  +3.31%  "

$ws.Range("D18").Value = "42.221.07"
$ws.Range("E18").Value = "  +1.29%  "

$ws.Range("D19").Value = "'12.51"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("D21").Value = "'5.97"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").Value = "'67.61"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").Value = "'240.64"

$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("D25").Value = "'1.95"
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'23.84"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").Value = "'37.88"
$ws.Range("E28").Value = "  +6.93%  "

$ws.Range("D29").Value = "'9.50"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").Value = "'2.11"
$ws.Range("E30").Value = "  +1.47%  "

$ws.Range("D31").Value = "'159.51"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").Value = "'5.23"
$ws.Range("E32").Value = "  -0.34%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  +3.82%  "

$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("D36").Value = "'17.01"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("D39").Value = "'1.85"
$ws.Range("E39").Value = "  +3.25%  "

$ws.Range("E40").Value = "  -1.30%  "

$ws.Range("E41").Value = "  +4.87%  "

$ws.Range("E42").Value = "  +13.84%  "

$ws.Range("D43").Value = "1.998.64"
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.09"
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0286"
$ws.Range("E45").Value = "  +2.31%  "

$ws.Range("D46").Value = "'2.99"
$ws.Range("E46").Value = "  +3.34%  "

$ws.Range("D47").Value = "'9.98"
$ws.Range("E47").Value = "  -3.48%  "

$ws.Range("D48").Value = "'53.05"
$ws.Range("E48").Value = "  +1.46%  "

$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("D50").Value = "'72.17"
$ws.Range("E50").Value = "  +0.17%  "

$ws.Range("D51").Value = "'91.92"
$ws.Range("E51").Value = "  +0.81%  "
